$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update recalculated KPI values (rerun of FA.py with a fix)
$ws.Range("AB2").Value = 8320.6954544990203
$ws.Range("AC2").Value = 8390.6257986399105
$ws.Range("AD2").Value = 8183.1445506885702
$ws.Range("AM2").Value = 1.12509774144034
$ws.Range("AN2").Value = 32.574527050575703
$ws.Range("AR2").Value = 56027.680839403802
$ws.Range("AS2").Value = 32336.716786871599
$ws.Range("AT2").Value = 3902.78445407138
$ws.Range("H2").Value = 8154.2815454090396
$ws.Range("I2").Value = 7400.5319544003996
$ws.Range("J2").Value = 6472.8673395946598
$ws.Range("T2").Value = 0.79100000000000004
$ws.Range("AA3").Value = 14841.5788166
$ws.Range("AB3").Value = 14298.9775580934
$ws.Range("AC3").Value = 15087.546777265799
$ws.Range("AD3").Value = 14149.017830876701
$ws.Range("AM3").Value = -0.0000000000509201443868273
$ws.Range("AN3").Value = 119.49138225761
$ws.Range("AR3").Value = 84590.764571542706
$ws.Range("AS3").Value = 57457.690723968197
$ws.Range("AT3").Value = 460.67465902816002
$ws.Range("F3").Value = 11326.4726206172
$ws.Range("G3").Value = 14663.4798708008
$ws.Range("H3").Value = 13012.069577864901
$ws.Range("I3").Value = 11376.010270058399
$ws.Range("J3").Value = 8899.7322322024793
$ws.Range("R3").Value = 0.91
$ws.Range("S3").Value = 0.754
$ws.Range("T3").Value = 0.629
$ws.Range("Z3").Value = 11947.7559289211
$ws.Range("AA4").Value = 23135.364620228302
$ws.Range("AB4").Value = 27407.792237613001
$ws.Range("AC4").Value = 25974.0714490725
$ws.Range("AD4").Value = 23828.828395672899
$ws.Range("AL4").Value = 29.939842693673999
$ws.Range("AM4").Value = 1.4722042540484901
$ws.Range("AN4").Value = 32.395950729590403
$ws.Range("AR4").Value = 85096.275590217905
$ws.Range("AS4").Value = 57777.725226495997
$ws.Range("AT4").Value = 275.19814288144102
$ws.Range("G4").Value = 14644.6858046045
$ws.Range("H4").Value = 27407.792237613001
$ws.Range("I4").Value = 23766.275375901299
$ws.Range("J4").Value = 19277.522172099401
$ws.Range("R4").Value = 1
$ws.Range("S4").Value = 0.91500000000000004
$ws.Range("T4").Value = 0.80900000000000005
$ws.Range("AB5").Value = 14175.075921731101
$ws.Range("AC5").Value = 42234.771077121499
$ws.Range("AD5").Value = 45459.172717627502
$ws.Range("AN5").Value = 11.9152096590337
$ws.Range("AR5").Value = 85437.141945819996
$ws.Range("AS5").Value = 58071.311870656202
$ws.Range("AT5").Value = 227.91843144007899
$ws.Range("H5").Value = 1261.5817570340701
$ws.Range("I5").Value = 38898.224162028899
$ws.Range("J5").Value = 45277.336026757002
$ws.Range("R5").Value = 0.088999999999999996
$ws.Range("AC6").Value = 13363.735290308499
$ws.Range("AD6").Value = 83877.858197528796
$ws.Range("AO6").Value = -0.0000000000099015242107961194
$ws.Range("AR6").Value = 85575.052579397801
$ws.Range("AS6").Value = 58208.288070174298
$ws.Range("AT6").Value = 226.983997380397
$ws.Range("I6").Value = 1697.19438186918
$ws.Range("J6").Value = 83877.858197528796
$ws.Range("AG7").Value = 8440.2726089314292
$ws.Range("AH7").Value = 8090.3085101509196
$ws.Range("AI7").Value = 7295.6082899705098
$ws.Range("AJ7").Value = 7772.4263526985396
$ws.Range("AU7").Value = 0.00105473173441434
$ws.Range("AV7").Value = 32.798839934460602
$ws.Range("AZ7").Value = 69928.603066785698
$ws.Range("BA7").Value = 48163.570247952499
$ws.Range("BB7").Value = 2093.4516831158498
$ws.Range("I7").Value = 8440.2726089314292
$ws.Range("J7").Value = 7669.6124676230702
$ws.Range("K7").Value = 6033.4680703964304
$ws.Range("L7").Value = 5285.2499198350197
$ws.Range("X7").Value = 0.68
$ws.Range("AD8").Value = 4378.8265096777304
$ws.Range("AE8").Value = 16999.999994273701
$ws.Range("AF8").Value = 16999.999999999902
$ws.Range("AG8").Value = 16313.7494493782
$ws.Range("AH8").Value = 14043.426643761401
$ws.Range("AI8").Value = 13974.432207784699
$ws.Range("AJ8").Value = 12346.748493152299
$ws.Range("AT8").Value = -0.000000000034327285982361501
$ws.Range("AU8").Value = -0.00000000013295010858928701
$ws.Range("AV8").Value = 46.489712965007598
$ws.Range("AZ8").Value = 85068.802309189399
$ws.Range("BA8").Value = 61589.603311129002
$ws.Range("BB8").Value = 379.28550388842598
$ws.Range("F8").Value = 4378.8265096777304
$ws.Range("G8").Value = 16999.999994273701
$ws.Range("H8").Value = 16999.999999999902
$ws.Range("I8").Value = 16313.7494493782
$ws.Range("J8").Value = 12849.7353790417
$ws.Range("K8").Value = 10019.667892981601
$ws.Range("L8").Value = 7506.8230838366098
$ws.Range("V8").Value = 0.91500000000000004
$ws.Range("W8").Value = 0.71699999999999997
$ws.Range("X8").Value = 0.60799999999999998
$ws.Range("AF9").Value = 17428.743817939099
$ws.Range("AG9").Value = 5976.5998846025796
$ws.Range("AH9").Value = 30626.2667252063
$ws.Range("AI9").Value = 26628.877416444899
$ws.Range("AJ9").Value = 23225.302098776701
$ws.Range("AU9").Value = -0.0000000000031338774934949201
$ws.Range("AV9").Value = 10.788178300779
$ws.Range("AW9").Value = -0.0000000000000248717004138803
$ws.Range("AZ9").Value = 80646.822443359095
$ws.Range("BA9").Value = 57222.045606976899
$ws.Range("BB9").Value = 433.707665566819
$ws.Range("H9").Value = 348.574876358783
$ws.Range("I9").Value = 5958.6700849487697
$ws.Range("J9").Value = 30565.014191755901
$ws.Range("K9").Value = 25217.5469133732
$ws.Range("L9").Value = 18557.016376922598
$ws.Range("W9").Value = 0.94699999999999995
$ws.Range("X9").Value = 0.79900000000000004
$ws.Range("AH10").Value = 28884.6721540818
$ws.Range("AI10").Value = 40033.662544765597
$ws.Range("AJ10").Value = 43545.316660796801
$ws.Range("AU10").Value = -0.0000000000048547826777381104
$ws.Range("AV10").Value = 10.667175412020301
$ws.Range("AZ10").Value = 80198.016691661993
$ws.Range("BA10").Value = 56738.276479994602
$ws.Range("BB10").Value = 398.74429028173699
$ws.Range("J10").Value = 1415.3489355500001
$ws.Range("K10").Value = 36630.801228460499
$ws.Range("L10").Value = 42151.866527651298
$ws.Range("AI11").Value = 28342.3053324295
$ws.Range("AJ11").Value = 77889.032402131794
$ws.Range("AV11").Value = 0.000000000088424506117368198
$ws.Range("AW11").Value = -0.000000000020527847732359801
$ws.Range("AX11").Value = -0.000000000000064184633888901396
$ws.Range("AZ11").Value = 80213.101439390899
$ws.Range("BA11").Value = 56738.259512596902
$ws.Range("BB11").Value = 383.64257515503601
$ws.Range("K11").Value = 2324.0690372592198
$ws.Range("L11").Value = 77889.032402131794

# Cells where Excel auto-applied scientific notation formatting for tiny residuals
$ws.Range("AM3").NumberFormat = "0.00E+00"
$ws.Range("AO6").NumberFormat = "0.00E+00"
$ws.Range("AT8").NumberFormat = "0.00E+00"
$ws.Range("AU8").NumberFormat = "0.00E+00"
$ws.Range("AU9").NumberFormat = "0.00E+00"
$ws.Range("AW9").NumberFormat = "0.00E+00"
$ws.Range("AU10").NumberFormat = "0.00E+00"
$ws.Range("AV11").NumberFormat = "0.00E+00"
$ws.Range("AW11").NumberFormat = "0.00E+00"
$ws.Range("AX11").NumberFormat = "0.00E+00"

# Update the active selection
$ws.Range("H17").Select()
